$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -2831
$ws.Range("N13").Value = -3338

$ws.Range("H41").Value = 1345.8
$ws.Range("I41").Value = 1080
$ws.Range("J41").Value = 1611.6
$ws.Range("K41").Value = 1080
$ws.Range("L41").Value = 1611.6
$ws.Range("M41").Value = -640
$ws.Range("N41").Value = -2491.6

$ws.Range("H51").Value = 5968.625
$ws.Range("I51").Value = 4999
$ws.Range("J51").Value = 6107.143
$ws.Range("K51").Value = 4999
$ws.Range("L51").Value = 6107.143
$ws.Range("M51").Value = -4515
$ws.Range("N51").Value = -7075.143

$ws.Range("H86").Value = 7572.222
$ws.Range("I86").Value = 8771.429
$ws.Range("J86").Value = 3375
$ws.Range("K86").Value = 8771.429
$ws.Range("L86").Value = 3375
$ws.Range("M86").Value = -7648.429

$ws.Range("H89").Value = 7572.222
$ws.Range("I89").Value = 8771.429
$ws.Range("J89").Value = 3375
$ws.Range("K89").Value = 43857.145
$ws.Range("L89").Value = 16875
$ws.Range("M89").Value = -38241.145

$ws.Range("H92").Value = 2489.7334
$ws.Range("I92").Value = 2642.3333
$ws.Range("J92").Value = 2388
$ws.Range("K92").Value = 2642.3333
$ws.Range("L92").Value = 2388
$ws.Range("M92").Value = -1394.3333
$ws.Range("N92").Value = -4884

$ws.Range("H113").Value = 6641.52
$ws.Range("I113").Value = 7479.294
$ws.Range("J113").Value = 4861.25
$ws.Range("K113").Value = 7479.294
$ws.Range("L113").Value = 4861.25
$ws.Range("M113").Value = -4225.294
$ws.Range("N113").Value = -11369.25

$ws.Range("H129").Value = 1298.6666
$ws.Range("I129").Value = 683.1667
$ws.Range("J129").Value = 1544.8667
$ws.Range("K129").Value = 2049.5001
$ws.Range("L129").Value = 4634.6001
$ws.Range("M129").Value = 2950.4999
$ws.Range("N129").Value = -14634.6001

$ws.Range("H136").Value = 85389.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 85389.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 85389.5
$ws.Range("N136").Value = -95589.5

$ws.Range("H137").Value = 10889
$ws.Range("I137").Value = 2772.3704
$ws.Range("J137").Value = 24585.812
$ws.Range("K137").Value = 8317.111199999999
$ws.Range("L137").Value = 73757.436
$ws.Range("M137").Value = -5767.111199999999
$ws.Range("N137").Value = -78857.436

$ws.Range("H138").Value = 3913.5454
$ws.Range("I138").Value = 4092.5
$ws.Range("J138").Value = 3811.2856
$ws.Range("K138").Value = 12277.5
$ws.Range("L138").Value = 11433.8568
$ws.Range("M138").Value = -7137.5
$ws.Range("N138").Value = -21713.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2994.7222
$ws.Range("I2").Value = 2738.6
$ws.Range("J2").Value = 4275.3335
$ws.Range("K2").Value = 2738.6
$ws.Range("L2").Value = 4275.3335
$ws.Range("M2").Value = -2625.6
$ws.Range("N2").Value = -4501.3335

$ws.Range("H8").Value = 626018.75
$ws.Range("I8").Value = 626018.75
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 626018.75
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -625874.75

$ws.Range("H32").Value = 4783.84
$ws.Range("I32").Value = 2076
$ws.Range("J32").Value = 8845.6
$ws.Range("K32").Value = 2076
$ws.Range("L32").Value = 8845.6
$ws.Range("M32").Value = -1789
$ws.Range("N32").Value = -9419.6

$ws.Range("H36").Value = 1666.3334
$ws.Range("I36").Value = 1666.3334
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1666.3334
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1320.3334

$ws.Range("H61").Value = 478662
$ws.Range("I61").Value = 2607.7693
$ws.Range("J61").Value = 2954144
$ws.Range("K61").Value = 2607.7693
$ws.Range("L61").Value = 2954144
$ws.Range("M61").Value = -2395.7693

$ws.Range("H74").Value = 7417.45
$ws.Range("I74").Value = 2554.5757
$ws.Range("J74").Value = 30342.428
$ws.Range("K74").Value = 2554.5757
$ws.Range("L74").Value = 30342.428
$ws.Range("M74").Value = -1680.5757

$ws.Range("H77").Value = 7417.45
$ws.Range("I77").Value = 2554.5757
$ws.Range("J77").Value = 30342.428
$ws.Range("K77").Value = 12772.8785
$ws.Range("L77").Value = 151712.14
$ws.Range("M77").Value = -8404.878499999999

$ws.Range("H116").Value = 2994.7222
$ws.Range("I116").Value = 2738.6
$ws.Range("J116").Value = 4275.3335
$ws.Range("K116").Value = 2738.6
$ws.Range("L116").Value = 4275.3335
$ws.Range("M116").Value = -444.5999999999999
$ws.Range("N116").Value = -8863.333500000001

$ws.Range("H132").Value = 650396.3
$ws.Range("I132").Value = 3082.9268
$ws.Range("J132").Value = 3599268.5
$ws.Range("K132").Value = 9248.7804
$ws.Range("L132").Value = 10797805.5
$ws.Range("M132").Value = -6718.7804
$ws.Range("N132").Value = -10802865.5

$ws.Range("H136").Value = 478662
$ws.Range("I136").Value = 2607.7693
$ws.Range("J136").Value = 2954144
$ws.Range("K136").Value = 7823.3079
$ws.Range("L136").Value = 8862432
$ws.Range("M136").Value = -5273.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2994.7222
$ws.Range("I3").Value = 2738.6
$ws.Range("J3").Value = 4275.3335
$ws.Range("K3").Value = 2738.6
$ws.Range("L3").Value = 4275.3335
$ws.Range("M3").Value = -2624.6
$ws.Range("N3").Value = -4503.3335

$ws.Range("H33").Value = 2512
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 2024
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 2024
$ws.Range("M33").Value = -2664
$ws.Range("N33").Value = -2696

$ws.Range("H105").Value = 3204.0908
$ws.Range("I105").Value = 2839.9412
$ws.Range("J105").Value = 4442.2
$ws.Range("K105").Value = 2839.9412
$ws.Range("L105").Value = 4442.2
$ws.Range("M105").Value = -1092.9412

$ws.Range("H134").Value = 11866.869
$ws.Range("I134").Value = 7807.143
$ws.Range("J134").Value = 24784.182
$ws.Range("K134").Value = 23421.429
$ws.Range("L134").Value = 74352.546
$ws.Range("M134").Value = -20886.429
$ws.Range("N134").Value = -79422.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 112442.91
$ws.Range("I31").Value = 202669.6
$ws.Range("J31").Value = 37254
$ws.Range("K31").Value = 202669.6
$ws.Range("L31").Value = 37254
$ws.Range("M31").Value = -202374.6

$ws.Range("H32").Value = 9750
$ws.Range("I32").Value = 9750
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9750
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9434

$ws.Range("H34").Value = 112442.91
$ws.Range("I34").Value = 202669.6
$ws.Range("J34").Value = 37254
$ws.Range("K34").Value = 202669.6
$ws.Range("L34").Value = 37254
$ws.Range("M34").Value = -202467.6

$ws.Range("H41").Value = 25000
$ws.Range("I41").Value = 25000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 25000
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -24572

$ws.Range("H58").Value = 9997.388999999999
$ws.Range("I58").Value = 5086.3076
$ws.Range("J58").Value = 22766.2
$ws.Range("K58").Value = 5086.3076
$ws.Range("L58").Value = 22766.2
$ws.Range("M58").Value = -4883.3076
$ws.Range("N58").Value = -23172.2

$ws.Range("H99").Value = 3580.4167
$ws.Range("I99").Value = 3583.5
$ws.Range("J99").Value = 3574.25
$ws.Range("K99").Value = 3583.5
$ws.Range("L99").Value = 3574.25
$ws.Range("M99").Value = -2085.5
$ws.Range("N99").Value = -6570.25

$ws.Range("H126").Value = 3580.4167
$ws.Range("I126").Value = 3583.5
$ws.Range("J126").Value = 3574.25
$ws.Range("K126").Value = 10750.5
$ws.Range("L126").Value = 10722.75
$ws.Range("M126").Value = -8280.5
$ws.Range("N126").Value = -15662.75

$ws.Range("H134").Value = 4572.5557
$ws.Range("I134").Value = 1841.1177
$ws.Range("J134").Value = 51007
$ws.Range("K134").Value = 5523.3531
$ws.Range("L134").Value = 153021
$ws.Range("M134").Value = -2988.3531
$ws.Range("N134").Value = -158091

$ws.Range("H136").Value = 9997.388999999999
$ws.Range("I136").Value = 5086.3076
$ws.Range("J136").Value = 22766.2
$ws.Range("K136").Value = 15258.9228
$ws.Range("L136").Value = 68298.60000000001
$ws.Range("M136").Value = -12708.9228
$ws.Range("N136").Value = -73398.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 285.57144
$ws.Range("I2").Value = 137.5
$ws.Range("J2").Value = 483
$ws.Range("K2").Value = 825
$ws.Range("L2").Value = 2898
$ws.Range("M2").Value = -712
$ws.Range("N2").Value = -3124

$ws.Range("H10").Value = 104.25
$ws.Range("I10").Value = 82.5
$ws.Range("J10").Value = 169.5
$ws.Range("K10").Value = 247.5
$ws.Range("L10").Value = 508.5
$ws.Range("M10").Value = -108.5
$ws.Range("N10").Value = -786.5

$ws.Range("H68").Value = 1152.0435
$ws.Range("I68").Value = 1028.8572
$ws.Range("J68").Value = 1205.9375
$ws.Range("K68").Value = 3086.5716
$ws.Range("L68").Value = 3617.8125
$ws.Range("M68").Value = -2275.5716
$ws.Range("N68").Value = -5239.8125

$ws.Range("H71").Value = 1152.0435
$ws.Range("I71").Value = 1028.8572
$ws.Range("J71").Value = 1205.9375
$ws.Range("K71").Value = 9259.7148
$ws.Range("L71").Value = 10853.4375
$ws.Range("M71").Value = -5203.7148
$ws.Range("N71").Value = -18965.4375

$ws.Range("H113").Value = 718.6667
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 628
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 1884
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -6224

$ws.Range("H122").Value = 13450730
$ws.Range("I122").Value = 20763602
$ws.Range("J122").Value = 4048466.8
$ws.Range("K122").Value = 186872418
$ws.Range("L122").Value = 36436201.2
$ws.Range("M122").Value = -186869968
$ws.Range("N122").Value = -36441101.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61.666668
$ws.Range("I2").Value = 55.2
$ws.Range("J2").Value = 69.75
$ws.Range("K2").Value = 55.2
$ws.Range("L2").Value = 69.75
$ws.Range("M2").Value = 57.8
$ws.Range("N2").Value = -295.75

$ws.Range("H29").Value = 8699.75
$ws.Range("I29").Value = 900
$ws.Range("J29").Value = 16499.5
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 16499.5
$ws.Range("M29").Value = -610
$ws.Range("N29").Value = -17079.5

$ws.Range("H132").Value = 611616.75
$ws.Range("I132").Value = 2453.5
$ws.Range("J132").Value = 1708110.6
$ws.Range("K132").Value = 7360.5
$ws.Range("L132").Value = 5124331.800000001
$ws.Range("M132").Value = -4830.5
$ws.Range("N132").Value = -5129391.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3412.6
$ws.Range("I40").Value = 2546.7144
$ws.Range("J40").Value = 5433
$ws.Range("K40").Value = 2546.7144
$ws.Range("L40").Value = 5433
$ws.Range("M40").Value = -2410.7144
$ws.Range("N40").Value = -5705

$ws.Range("H93").Value = 16702.5
$ws.Range("I93").Value = 16702.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 16702.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -15454.5

$ws.Range("H122").Value = 6686
$ws.Range("I122").Value = 7472.4546
$ws.Range("J122").Value = 3802.3333
$ws.Range("K122").Value = 22417.3638
$ws.Range("L122").Value = 11406.9999
$ws.Range("M122").Value = -19967.3638

$ws.Range("H132").Value = 1426878.8
$ws.Range("I132").Value = 2899.6365
$ws.Range("J132").Value = 2993255.8
$ws.Range("K132").Value = 8698.9095
$ws.Range("L132").Value = 8979767.399999999
$ws.Range("M132").Value = -6168.9095
$ws.Range("N132").Value = -8984827.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 232591.67
$ws.Range("I136").Value = 1261.2693
$ws.Range("J136").Value = 533321.2
$ws.Range("K136").Value = 3783.8079
$ws.Range("L136").Value = 1599963.6
$ws.Range("M136").Value = -1233.8079
